# Scheduled runner update: refresh market-price/profit columns (H-N) on
# several leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with
# newly pulled values. Some rows gain/lose the profit (M/N) columns
# entirely depending on whether a meaningful profit could be computed.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2912.4285
$ws.Range("I107").Value = 2912.4285
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2912.4285
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -992.4285
$ws.Range("N107").ClearContents()
$ws.Range("H112").Value = 4804.2812
$ws.Range("J112").Value = 5181.3105
$ws.Range("L112").Value = 15543.9315
$ws.Range("N112").Value = -17759.9315

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H32").Value = 4046.2673
$ws.Range("I32").Value = 2495.257
$ws.Range("J32").Value = 10831.9375
$ws.Range("K32").Value = 2495.257
$ws.Range("L32").Value = 10831.9375
$ws.Range("M32").Value = -2208.257
$ws.Range("N32").Value = -11405.9375
$ws.Range("H44").Value = 6944
$ws.Range("H45").Value = 76228.86
$ws.Range("I45").Value = 112693
$ws.Range("K45").Value = 112693
$ws.Range("M45").Value = -112316
$ws.Range("H55").Value = 37524
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H63").Value = 3262.2727
$ws.Range("I63").Value = 3198.5
$ws.Range("J63").Value = 3900
$ws.Range("K63").Value = 3198.5
$ws.Range("L63").Value = 3900
$ws.Range("M63").Value = -2512.5
$ws.Range("N63").Value = -5272
$ws.Range("H66").Value = 3262.2727
$ws.Range("I66").Value = 3198.5
$ws.Range("J66").Value = 3900
$ws.Range("K66").Value = 15992.5
$ws.Range("L66").Value = 19500
$ws.Range("M66").Value = -12560.5
$ws.Range("N66").Value = -26364
$ws.Range("H74").Value = 190476.9
$ws.Range("I74").Value = 167323.83
$ws.Range("J74").Value = 225206.5
$ws.Range("K74").Value = 167323.83
$ws.Range("L74").Value = 225206.5
$ws.Range("M74").Value = -166449.83
$ws.Range("N74").Value = -226954.5
$ws.Range("H77").Value = 190476.9
$ws.Range("I77").Value = 167323.83
$ws.Range("J77").Value = 225206.5
$ws.Range("K77").Value = 836619.1499999999
$ws.Range("L77").Value = 1126032.5
$ws.Range("M77").Value = -832251.1499999999
$ws.Range("N77").Value = -1134768.5
$ws.Range("H97").Value = 15703.685
$ws.Range("I97").Value = 14477.692
$ws.Range("J97").Value = 18360
$ws.Range("K97").Value = 14477.692
$ws.Range("L97").Value = 18360
$ws.Range("M97").Value = -13981.692
$ws.Range("N97").Value = -19352
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H119").Value = 88599.44500000001
$ws.Range("J119").Value = 88599.44500000001
$ws.Range("L119").Value = 88599.44500000001
$ws.Range("N119").Value = -98275.44500000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 95000
$ws.Range("J110").Value = 95000
$ws.Range("L110").Value = 95000
$ws.Range("N110").Value = -103180

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 16209.8
$ws.Range("J41").Value = 32495
$ws.Range("L41").Value = 32495
$ws.Range("N41").Value = -33351
$ws.Range("H50").Value = 2621.5
$ws.Range("J50").Value = 2621.5
$ws.Range("L50").Value = 2621.5
$ws.Range("N50").Value = -3871.5
$ws.Range("H51").Value = 29798.666
$ws.Range("J51").Value = 29798.666
$ws.Range("L51").Value = 29798.666
$ws.Range("N51").Value = -31270.666
$ws.Range("H58").Value = 3996.6667
$ws.Range("I58").Value = 3497.5
$ws.Range("J58").Value = 4995
$ws.Range("K58").Value = 3497.5
$ws.Range("L58").Value = 4995
$ws.Range("M58").Value = -3294.5
$ws.Range("N58").Value = -5401
$ws.Range("H59").Value = 37399.9
$ws.Range("I59").Value = 50000
$ws.Range("J59").Value = 35999.89
$ws.Range("K59").Value = 50000
$ws.Range("L59").Value = 35999.89
$ws.Range("M59").Value = -48855
$ws.Range("N59").Value = -38289.89
$ws.Range("H61").Value = 29798.666
$ws.Range("J61").Value = 29798.666
$ws.Range("L61").Value = 29798.666
$ws.Range("N61").Value = -30494.666
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H99").Value = 4247.5884
$ws.Range("I99").Value = 3703.5454
$ws.Range("K99").Value = 3703.5454
$ws.Range("M99").Value = -2205.5454
$ws.Range("H107").Value = 2564.0715
$ws.Range("I107").Value = 2737.2222
$ws.Range("K107").Value = 2737.2222
$ws.Range("M107").Value = -817.2222000000002
$ws.Range("H126").Value = 4247.5884
$ws.Range("I126").Value = 3703.5454
$ws.Range("K126").Value = 11110.6362
$ws.Range("M126").Value = -8640.636200000001
$ws.Range("H134").Value = 26727.475
$ws.Range("I134").Value = 36952.19
$ws.Range("J134").Value = 4573.9165
$ws.Range("K134").Value = 110856.57
$ws.Range("L134").Value = 13721.7495
$ws.Range("M134").Value = -108321.57
$ws.Range("N134").Value = -18791.7495
$ws.Range("H136").Value = 3996.6667
$ws.Range("I136").Value = 3497.5
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 10492.5
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -7942.5
$ws.Range("N136").Value = -20085

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1246.6666
$ws.Range("I14").Value = 1246.6666
$ws.Range("K14").Value = 3739.9998
$ws.Range("M14").Value = -3566.9998
$ws.Range("H132").Value = 2965.5454
$ws.Range("I132").Value = 2820.375
$ws.Range("K132").Value = 25383.375
$ws.Range("M132").Value = -22853.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 514997.5
$ws.Range("I62").Value = 514997.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 514997.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -514311.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 514997.5
$ws.Range("I65").Value = 514997.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1544992.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -1541560.5
$ws.Range("N65").ClearContents()
$ws.Range("H80").Value = 2682.6
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2758.4443
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2758.4443
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4754.4443
$ws.Range("H83").Value = 2682.6
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2758.4443
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 13792.2215
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -23776.2215

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 10009
$ws.Range("I4").Value = 10009
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 10009
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -9896
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 4125.2666
$ws.Range("I7").Value = 2234.6365
$ws.Range("K7").Value = 2234.6365
$ws.Range("M7").Value = -2122.6365
$ws.Range("H28").Value = 10009
$ws.Range("I28").Value = 10009
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 10009
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -9777
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 10009
$ws.Range("I37").Value = 10009
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 10009
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -9902
$ws.Range("N37").ClearContents()
$ws.Range("H126").Value = 4125.2666
$ws.Range("I126").Value = 2234.6365
$ws.Range("K126").Value = 6703.9095
$ws.Range("M126").Value = -4233.9095
$ws.Range("H136").Value = 405600.6
$ws.Range("I136").Value = 405600.6
$ws.Range("K136").Value = 1216801.8
$ws.Range("M136").Value = -1214251.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1128.7307
$ws.Range("I113").Value = 818.5
$ws.Range("K113").Value = 2455.5
$ws.Range("M113").Value = -285.5
$ws.Range("H119").Value = 69739.39999999999
$ws.Range("J119").Value = 69739.39999999999
$ws.Range("L119").Value = 69739.39999999999
$ws.Range("N119").Value = -79415.39999999999
